$d = $word.ActiveDocument

$replacements = @(
    @{old="61×53=3233"; new="27×65=1755"},
    @{old="83×50=4150"; new="98×11=1078"},
    @{old="84×43=3612"; new="37×47=1739"},
    @{old="94×73=6862"; new="88×21=1848"},
    @{old="71×72=5112"; new="91×27=2457"},
    @{old="68×49=3332"; new="40×30=1200"},
    @{old="60×81=4860"; new="79×63=4977"},
    @{old="62×22=1364"; new="84×97=8148"},
    @{old="65×86=5590"; new="57×28=1596"},
    @{old="94×43=4042"; new="94×87=8178"},
    @{old="92×68=6256"; new="65×18=1170"},
    @{old="54×32=1728"; new="20×72=1440"},
    @{old="27×87=2349"; new="73×88=6424"},
    @{old="39×44=1716"; new="59×93=5487"},
    @{old="96×66=6336"; new="99×25=2475"},
    @{old="55×22=1210"; new="78×47=3666"},
    @{old="16×33=528";  new="30×98=2940"},
    @{old="75×31=2325"; new="41×76=3116"},
    @{old="74×89=6586"; new="33×62=2046"},
    @{old="83×69=5727"; new="32×90=2880"},
    @{old="56×34=1904"; new="65×41=2665"},
    @{old="23×64=1472"; new="89×86=7654"},
    @{old="79×95=7505"; new="71×25=1775"},
    @{old="45×17=765";  new="65×94=6110"},
    @{old="48×19=912";  new="94×85=7990"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
